# Product Backlog: fill in real user-story rows (Online Ordering Interface
# epic) in place of the single "Feature 1" placeholder row, extending the
# table from row 4 down through row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend formatting down to the new rows -----------------------------
# Rows 7:15 already carry the canonical "data row" look (thin top border on
# A, thin border box on B:E, centered border on F). Rows 5:6 are left over
# from an earlier draft and use a slightly different style, and rows 16:20
# do not exist yet. Write placeholders so the sheet's used range covers
# A5:F20, then copy the row-7 formatting down across every data row so the
# whole block (5:20) is visually consistent before we fill in real values.
$ws.Range("A16:F20").Value = ""
$ws.Range("A7:F7").Copy()
$ws.Range("A5:F20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Feature / story rows -------------------------------------------------
$rows = @(
    @{ Row = 4;  Name = "Online Ordering Interface";                      Priority = "High";   Points = 8 },
    @{ Row = 5;  Name = "Develop customer-facing web/API platform";       Priority = "High";   Points = 5 },
    @{ Row = 6;  Name = "Enable guest checkout for new users";            Priority = "High";   Points = 3 },
    @{ Row = 7;  Name = "Design menu display with dish metadata";         Priority = "Medium"; Points = 5 },
    @{ Row = 8;  Name = "Provide takeout and delivery selection options"; Priority = "High";   Points = 3 },
    @{ Row = 9;  Name = "Integrate real-time order tracking system";      Priority = "High";   Points = 8 },
    @{ Row = 10; Name = "Setup multi-method payment processing";          Priority = "High";   Points = 8 },
    @{ Row = 11; Name = "Create role-based dashboards (staff/chef/admin)";Priority = "High";   Points = 8 },
    @{ Row = 12; Name = "Build real-time kitchen order display system";   Priority = "High";   Points = 5 },
    @{ Row = 13; Name = "Develop admin analytics dashboard";              Priority = "Medium"; Points = 5 },
    @{ Row = 14; Name = "Add customer review/rating module";              Priority = "Medium"; Points = 3 },
    @{ Row = 15; Name = "Implement feedback reporting and analysis tools";Priority = "Medium"; Points = 5 },
    @{ Row = 16; Name = "Enable coupon and promo code engine";            Priority = "Low";    Points = 3 },
    @{ Row = 17; Name = "Build promotion management dashboard";           Priority = "Medium"; Points = 5 },
    @{ Row = 18; Name = "Setup data analytics and sales reporting tools"; Priority = "Medium"; Points = 8 },
    @{ Row = 19; Name = "Create integrated training resources for staff"; Priority = "Low";    Points = 2 },
    @{ Row = 20; Name = "Provide comprehensive system documentation repo";Priority = "Medium"; Points = 3 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Name
    $ws.Range("B$row").Value = $r.Priority
    $ws.Range("C$row").Value = $r.Points
    $ws.Range("D$row").Value = "Yes"
    $ws.Range("E$row").Value = "Yes"
    $ws.Range("F$row").Value = 1
}

# --- 3. Data validation now spans the whole filled block (rows 4-20) -------
$ws.Cells.Validation.Delete()
$ws.Range("B4:B20").Validation.Add(3, 1, 1, '"High,Medium,Low"')
$ws.Range("D4:E20").Validation.Add(3, 1, 1, '"Yes,No"')

# --- 4. Selection cosmetics matching the authored workbook ------------------
$ws.Range("D26").Select()
